$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.427.32"
$ws.Range("E2").Value = "  +1.13%  "
$ws.Range("D3").Value = "'1.674.26"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("D4").Value = "'1.011"
$ws.Range("E4").Value = "  +0.88%  "
$ws.Range("D5").Value = "'221.26"
$ws.Range("E5").Value = "  +1.52%  "
$ws.Range("D6").Value = "'0.5353"
$ws.Range("E6").Value = "  +1.08%  "
$ws.Range("D8").Value = "'0.2676"
$ws.Range("E8").Value = "  +2.37%  "
$ws.Range("D9").Value = "'0.06414"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "'21.00"
$ws.Range("E10").Value = "  +2.83%  "
$ws.Range("D11").Value = "'0.07856"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "'4.550"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "'1.676.12"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("D14").Value = "'1.903.98"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "'0.5662"
$ws.Range("E15").Value = "  +3.17%  "
$ws.Range("D16").Value = "'0.0₅8207"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").Value = "'66.48"
$ws.Range("E17").Value = "  +1.64%  "
$ws.Range("D18").Value = "'26.472.38"
$ws.Range("E18").Value = "  +1.32%  "
$ws.Range("E19").Value = "  +0.82%  "
$ws.Range("D20").Value = "'4.729"
$ws.Range("E20").Value = "  +2.82%  "
$ws.Range("D21").Value = "'198.82"
$ws.Range("E21").Value = "  +4.02%  "
$ws.Range("E22").Value = "  +2.81%  "
$ws.Range("D23").Value = "'6.088"
$ws.Range("E23").Value = "  +0.99%  "
$ws.Range("D24").Value = "'1.012"
$ws.Range("E24").Value = "  +0.84%  "
$ws.Range("D25").Value = "'146.79"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'0.1234"
$ws.Range("E26").Value = "  +0.51%  "
$ws.Range("D27").Value = "'7.276"
$ws.Range("E27").Value = "  +0.81%  "
$ws.Range("D28").Value = "'16.25"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "'1.509"
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("D30").Value = "'0.05897"
$ws.Range("E30").Value = "  +2.16%  "
$ws.Range("D31").Value = "'1.289"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "'3.588"
$ws.Range("E32").Value = "  +0.79%  "
$ws.Range("D33").Value = "'3.318"
$ws.Range("E33").Value = "  +1.42%  "
$ws.Range("D34").Value = "'1.619"
$ws.Range("E34").Value = "  +0.91%  "
$ws.Range("D35").Value = "'0.9711"
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("D36").Value = "'2.851"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").Value = "'2.439"
$ws.Range("E37").Value = "  +1.10%  "
$ws.Range("D38").Value = "'0.5839"
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("D39").Value = "'0.01617"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").Value = "'1.079.25"
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("D41").Value = "'5.924"
$ws.Range("E41").Value = "  +2.43%  "
$ws.Range("D42").Value = "'0.8669"
$ws.Range("E42").Value = "  +1.28%  "
$ws.Range("E43").Value = "  +0.83%  "
$ws.Range("D44").Value = "'104.30"
$ws.Range("E44").Value = "  -0.39%  "
$ws.Range("D45").Value = "'1.813.52"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("D46").Value = "'58.56"
$ws.Range("E46").Value = "  +2.93%  "
$ws.Range("D47").Value = "'0.0₈107"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").Value = "'1.014"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("D49").Value = "'0.4414"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").Value = "'8.077"
$ws.Range("E50").Value = "  +2.45%  "

$ws.Range("D2:E51").ClearFormats()

